# Auto-generated PowerShell COM-interop script
$wb = $excel.ActiveWorkbook

# ---------- Sheet1: Student Summary ----------
$ws1 = $wb.Worksheets.Item("Student Summary")

# Copy formatting of row 10 down into the (currently empty) rows 11 and 12
$ws1.Range("A10:D10").Copy()
$ws1.Range("A11:D12").PasteSpecial(-4122)

# New "Course Code" / "DSPE603" / subject-code placeholder row
$ws1.Cells.Item(11,2).Value = "Course Code:"
$ws1.Cells.Item(11,3).Value = "DSPE603"
$ws1.Cells.Item(11,4).Value = "<---- Type Subject Code"

# New "Total Marks" / 40 / "Mid-Test Mark" row
$ws1.Cells.Item(12,2).Value = "Total Marks"
$ws1.Cells.Item(12,3).Value = 40
$ws1.Cells.Item(12,4).Value = "Mid-Test Mark"

# Round the Average Marks value
$ws1.Range("B17").Value = 27.24

# ---------- Sheet2: Slow Learners ----------
$ws2 = $wb.Worksheets.Item("Slow Learners")
$ws2.Cells.Item(2,1).Value = 28
$ws2.Cells.Item(2,2).Value = 2136110029
$ws2.Cells.Item(2,3).Value = "Balaganapathi A"
$ws2.Cells.Item(2,4).Value = 9
$ws2.Cells.Item(3,1).Value = 41
$ws2.Cells.Item(3,2).Value = 2136110042
$ws2.Cells.Item(3,3).Value = "Rajadurai P"
$ws2.Cells.Item(3,4).Value = 11
$ws2.Cells.Item(4,1).Value = 51
$ws2.Cells.Item(4,2).Value = 2236150003
$ws2.Cells.Item(4,3).Value = "Krishnakumar S"
$ws2.Cells.Item(4,4).Value = 11
$ws2.Cells.Item(5,1).Value = 17
$ws2.Cells.Item(5,2).Value = 2136110018
$ws2.Cells.Item(5,3).Value = "Sikanthkumar C"
$ws2.Cells.Item(5,4).Value = 13
$ws2.Cells.Item(6,1).Value = 34
$ws2.Cells.Item(6,2).Value = 2136110035
$ws2.Cells.Item(6,3).Value = "Manuneethi S"
$ws2.Cells.Item(6,4).Value = 14
$ws2.Cells.Item(7,1).Value = 50
$ws2.Cells.Item(7,2).Value = 2236150002
$ws2.Cells.Item(7,3).Value = "Mohammed Azees M"
$ws2.Cells.Item(7,4).Value = 14
# Remove the last (7th) slow-learner row -- that student no longer qualifies
$ws2.Rows.Item(8).Delete()

# ---------- Sheet3: Fast Learners ----------
$ws3 = $wb.Worksheets.Item("Fast Learners")
$ws3.Cells.Item(2,1).Value = 8
$ws3.Cells.Item(2,2).Value = 2136110008
$ws3.Cells.Item(2,3).Value = "Jananika B"
$ws3.Cells.Item(2,4).Value = 38
$ws3.Cells.Item(3,1).Value = 20
$ws3.Cells.Item(3,2).Value = 2136110021
$ws3.Cells.Item(3,3).Value = "Subhashini S"
$ws3.Cells.Item(3,4).Value = 38
$ws3.Cells.Item(4,1).Value = 9
$ws3.Cells.Item(4,2).Value = 2136110009
$ws3.Cells.Item(4,3).Value = "Kalaivani S"
$ws3.Cells.Item(4,4).Value = 38
$ws3.Cells.Item(5,1).Value = 1
$ws3.Cells.Item(5,2).Value = 2136110001
$ws3.Cells.Item(5,3).Value = "Aravind S"
$ws3.Cells.Item(5,4).Value = 37
$ws3.Cells.Item(6,1).Value = 21
$ws3.Cells.Item(6,2).Value = 2136110022
$ws3.Cells.Item(6,3).Value = "Suji Shri B"
$ws3.Cells.Item(6,4).Value = 37
$ws3.Cells.Item(7,1).Value = 4
$ws3.Cells.Item(7,2).Value = 2136110004
$ws3.Cells.Item(7,3).Value = "Deepakragavan J"
$ws3.Cells.Item(7,4).Value = 36
$ws3.Cells.Item(8,1).Value = 13
$ws3.Cells.Item(8,2).Value = 2136110013
$ws3.Cells.Item(8,3).Value = "Naveena A"
$ws3.Cells.Item(8,4).Value = 36
$ws3.Cells.Item(9,1).Value = 29
$ws3.Cells.Item(9,2).Value = 2136110030
$ws3.Cells.Item(9,3).Value = "Brijesh A"
$ws3.Cells.Item(9,4).Value = 35
$ws3.Cells.Item(10,1).Value = 47
$ws3.Cells.Item(10,2).Value = 2136110048
$ws3.Cells.Item(10,3).Value = "Bhuvanadurai M"
$ws3.Cells.Item(10,4).Value = 35
$ws3.Cells.Item(11,1).Value = 30
$ws3.Cells.Item(11,2).Value = 2136110031
$ws3.Cells.Item(11,3).Value = "Hitesh Kumar K A"
$ws3.Cells.Item(11,4).Value = 35
$ws3.Cells.Item(12,1).Value = 44
$ws3.Cells.Item(12,2).Value = 2136110045
$ws3.Cells.Item(12,3).Value = "Varsha V"
$ws3.Cells.Item(12,4).Value = 34
$ws3.Cells.Item(13,1).Value = 39
$ws3.Cells.Item(13,2).Value = 2136110040
$ws3.Cells.Item(13,3).Value = "Preethiga S"
$ws3.Cells.Item(13,4).Value = 34
$ws3.Cells.Item(14,1).Value = 19
$ws3.Cells.Item(14,2).Value = 2136110020
$ws3.Cells.Item(14,3).Value = "Srija D"
$ws3.Cells.Item(14,4).Value = 33
$ws3.Cells.Item(15,1).Value = 33
$ws3.Cells.Item(15,2).Value = 2136110034
$ws3.Cells.Item(15,3).Value = "Kaviraj M"
$ws3.Cells.Item(15,4).Value = 33
$ws3.Cells.Item(16,1).Value = 7
$ws3.Cells.Item(16,2).Value = 2136110007
$ws3.Cells.Item(16,3).Value = "Guruprasath V"
$ws3.Cells.Item(16,4).Value = 33
$ws3.Cells.Item(17,1).Value = 40
$ws3.Cells.Item(17,2).Value = 2136110041
$ws3.Cells.Item(17,3).Value = "Ragothaman R"
$ws3.Cells.Item(17,4).Value = 33
$ws3.Cells.Item(18,1).Value = 46
$ws3.Cells.Item(18,2).Value = 2136110047
$ws3.Cells.Item(18,3).Value = "Gowtham R"
$ws3.Cells.Item(18,4).Value = 33
$ws3.Cells.Item(19,1).Value = 48
$ws3.Cells.Item(19,2).Value = 2136110049
$ws3.Cells.Item(19,3).Value = "Kailashwaran R"
$ws3.Cells.Item(19,4).Value = 32
$ws3.Cells.Item(20,1).Value = 35
$ws3.Cells.Item(20,2).Value = 2136110036
$ws3.Cells.Item(20,3).Value = "Mohamed Suhail J"
$ws3.Cells.Item(20,4).Value = 32
$ws3.Cells.Item(21,1).Value = 16
$ws3.Cells.Item(21,2).Value = 2136110016
$ws3.Cells.Item(21,3).Value = "Nithya Sri R"
$ws3.Cells.Item(21,4).Value = 31
$ws3.Cells.Item(22,1).Value = 2
$ws3.Cells.Item(22,2).Value = 2136110002
$ws3.Cells.Item(22,3).Value = "Arulselvam C"
$ws3.Cells.Item(22,4).Value = 31
$ws3.Cells.Item(23,1).Value = 22
$ws3.Cells.Item(23,2).Value = 2136110023
$ws3.Cells.Item(23,3).Value = "Surya Prakash R"
$ws3.Cells.Item(23,4).Value = 31
$ws3.Cells.Item(24,1).Value = 14
$ws3.Cells.Item(24,2).Value = 2136110014
$ws3.Cells.Item(24,3).Value = "Nawin B"
$ws3.Cells.Item(24,4).Value = 31
$ws3.Cells.Item(25,1).Value = 18
$ws3.Cells.Item(25,2).Value = 2136110019
$ws3.Cells.Item(25,3).Value = "Sivaa Ganesh S"
$ws3.Cells.Item(25,4).Value = 31
$ws3.Cells.Item(26,1).Value = 49
$ws3.Cells.Item(26,2).Value = 2236150001
$ws3.Cells.Item(26,3).Value = "Dhanush B"
$ws3.Cells.Item(26,4).Value = 30
$ws3.Cells.Item(27,1).Value = 10
$ws3.Cells.Item(27,2).Value = 2136110010
$ws3.Cells.Item(27,3).Value = "Krishnapriya K"
$ws3.Cells.Item(27,4).Value = 29
$ws3.Cells.Item(28,1).Value = 38
$ws3.Cells.Item(28,2).Value = 2136110039
$ws3.Cells.Item(28,3).Value = "Pranav Varshan A T"
$ws3.Cells.Item(28,4).Value = 28
$ws3.Cells.Item(29,1).Value = 37
$ws3.Cells.Item(29,2).Value = 2136110038
$ws3.Cells.Item(29,3).Value = "Pradeep M"
$ws3.Cells.Item(29,4).Value = 28
$ws3.Cells.Item(30,1).Value = 3
$ws3.Cells.Item(30,2).Value = 2136110003
$ws3.Cells.Item(30,3).Value = "Ashik Jenly V L"
$ws3.Cells.Item(30,4).Value = 28
$ws3.Cells.Item(31,1).Value = 12
$ws3.Cells.Item(31,2).Value = 2136110012
$ws3.Cells.Item(31,3).Value = "Mohanraj D"
$ws3.Cells.Item(31,4).Value = 28
$ws3.Cells.Item(32,1).Value = 36
$ws3.Cells.Item(32,2).Value = 2136110037
$ws3.Cells.Item(32,3).Value = "Mohamed Tharif B"
$ws3.Cells.Item(32,4).Value = 27
$ws3.Cells.Item(33,1).Value = 5
$ws3.Cells.Item(33,2).Value = 2136110005
$ws3.Cells.Item(33,3).Value = "Devadharshini A"
$ws3.Cells.Item(33,4).Value = 27
$ws3.Cells.Item(34,1).Value = 6
$ws3.Cells.Item(34,2).Value = 2136110006
$ws3.Cells.Item(34,3).Value = "Dhinakaran R"
$ws3.Cells.Item(34,4).Value = 27
$ws3.Cells.Item(35,1).Value = 45
$ws3.Cells.Item(35,2).Value = 2136110046
$ws3.Cells.Item(35,3).Value = "Marikannan P"
$ws3.Cells.Item(35,4).Value = 26

